# Generate Report for Handback
#
# The "473ad77f-a431-4757-9ad7-6029ad044372" file has been handed back
# (translations delivered). Update the Overview + per-locale status
# sheets to reflect:
#   - Status: "Ready for handoff" -> "Handed back: in sync with en-US"
#   - Latest Target File / Latest Handback File now populated (with
#     hyperlinks to the source .md / translated .xlf respectively)
#   - Latest Handback DateTime populated (was the zero-date sentinel)

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: both the zh-cn and de-de status columns for this file
# flip to the new status text.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus

$zhcn.Hyperlinks.Add(
    $zhcn.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/926021af2663552d43e929f14597575b10aee4ae/e2e/473ad77f-a431-4757-9ad7-6029ad044372.md",
    $null,
    $null,
    "473ad77f-a431-4757-9ad7-6029ad044372.md"
) | Out-Null

$zhcn.Hyperlinks.Add(
    $zhcn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d82bcd176d6b523c11fb3348f62d8b45c0280020/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/473ad77f-a431-4757-9ad7-6029ad044372.b1493079a41cf13852acac0189b2bc9372264d57.zh-cn.xlf",
    $null,
    $null,
    "473ad77f-a431-4757-9ad7-6029ad044372.b1493079a41cf13852acac0189b2bc9372264d57.zh-cn.xlf"
) | Out-Null

$zhcn.Range("H2").Value = "2016-03-23 12:41:04"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus

$dede.Hyperlinks.Add(
    $dede.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/926021af2663552d43e929f14597575b10aee4ae/e2e/473ad77f-a431-4757-9ad7-6029ad044372.md",
    $null,
    $null,
    "473ad77f-a431-4757-9ad7-6029ad044372.md"
) | Out-Null

$dede.Hyperlinks.Add(
    $dede.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/360706b4da052ce039fa8ee26a45a2fcc8f804d6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/473ad77f-a431-4757-9ad7-6029ad044372.b1493079a41cf13852acac0189b2bc9372264d57.de-de.xlf",
    $null,
    $null,
    "473ad77f-a431-4757-9ad7-6029ad044372.b1493079a41cf13852acac0189b2bc9372264d57.de-de.xlf"
) | Out-Null

$dede.Range("H2").Value = "2016-03-23 12:41:10"
